# CryCompanywiseStockReport_1.xlsx update
#
# Refresh the stock-quantity figures on the report: the "Qty" column (F) and
# the dependent "Value" column (G = Rate(D) x Qty(F)) for a number of items
# were recalculated against newer stock counts. Every "Sub Total:" row's
# total (column B) and the final "Sub Total:"/"Grand Total:" rows were
# re-summed to reflect the updated item values. Two rows (283/284, the two
# "HUL-Bru Inst Poly 50g" batches) additionally had their batch no./rate/MRP
# values corrected along with quantity & value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F11").Value = 7
$ws.Range("G11").Value = 1525.51
$ws.Range("B12").Value = 5382.79
$ws.Range("F25").Value = 33
$ws.Range("G25").Value = 1496.22
$ws.Range("B27").Value = 10758.39
$ws.Range("F31").Value = 71
$ws.Range("G31").Value = 2217.33
$ws.Range("F36").Value = 24
$ws.Range("G36").Value = 4629.36
$ws.Range("F42").Value = 55
$ws.Range("G42").Value = 1689.6
$ws.Range("F45").Value = 72
$ws.Range("G45").Value = 6734.88
$ws.Range("F46").Value = 89
$ws.Range("G46").Value = 1460.49
$ws.Range("B56").Value = 41194.4
$ws.Range("F90").Value = 36
$ws.Range("G90").Value = 4144.32
$ws.Range("F92").Value = 0
$ws.Range("G92").Value = 0
$ws.Range("F104").Value = 150
$ws.Range("G104").Value = 15369
$ws.Range("F114").Value = 201
$ws.Range("G114").Value = 3911.46
$ws.Range("B115").Value = 251049.94
$ws.Range("F129").Value = 49
$ws.Range("G129").Value = 5103.35
$ws.Range("B133").Value = 12160.75
$ws.Range("F146").Value = 142
$ws.Range("G146").Value = 2794.56
$ws.Range("B151").Value = 23963.09
$ws.Range("F172").Value = 5
$ws.Range("G172").Value = 556.65
$ws.Range("F176").Value = 37
$ws.Range("G176").Value = 1461.13
$ws.Range("B185").Value = 17982.85
$ws.Range("F189").Value = 133
$ws.Range("G189").Value = 4126.99
$ws.Range("F194").Value = 27
$ws.Range("G194").Value = 2194.56
$ws.Range("F195").Value = 1
$ws.Range("G195").Value = 78.26000000000001
$ws.Range("B197").Value = 14110.87
$ws.Range("F214").Value = 37
$ws.Range("G214").Value = 3138.71
$ws.Range("F215").Value = 92
$ws.Range("G215").Value = 6339.72
$ws.Range("B217").Value = 10052.18
$ws.Range("F232").Value = 2
$ws.Range("G232").Value = 953.5
$ws.Range("B238").Value = 7701.62
$ws.Range("F261").Value = 41
$ws.Range("G261").Value = 1996.29
$ws.Range("F265").Value = 9
$ws.Range("G265").Value = 2449.08
$ws.Range("F273").Value = 2
$ws.Range("G273").Value = 68.92
$ws.Range("F276").Value = 2
$ws.Range("G276").Value = 934.38
$ws.Range("B279").Value = 117321.33
$ws.Range("B283").Value = 61610
$ws.Range("D283").Value = 102.71
$ws.Range("E283").Value = 122.71
$ws.Range("F283").Value = 199
$ws.Range("G283").Value = 20439.29
$ws.Range("B284").Value = 57077
$ws.Range("D284").Value = 93.08
$ws.Range("E284").Value = 111.2
$ws.Range("F284").Value = 1
$ws.Range("G284").Value = 93.08
$ws.Range("F288").Value = 24
$ws.Range("G288").Value = 1990.56
$ws.Range("F296").Value = 95
$ws.Range("G296").Value = 13021.65
$ws.Range("F312").Value = 158
$ws.Range("G312").Value = 17689.68
$ws.Range("F318").Value = 2
$ws.Range("G318").Value = 570.8200000000001
$ws.Range("F320").Value = 20
$ws.Range("G320").Value = 1776
$ws.Range("F323").Value = 133
$ws.Range("G323").Value = 13446.3
$ws.Range("F329").Value = 41
$ws.Range("G329").Value = 4863.83
$ws.Range("F330").Value = 93
$ws.Range("G330").Value = 5499.09
$ws.Range("B349").Value = 362059.11
$ws.Range("F386").Value = 3
$ws.Range("G386").Value = 2265
$ws.Range("B389").Value = 10490.52
$ws.Range("F402").Value = 47
$ws.Range("G402").Value = 4540.2
$ws.Range("F403").Value = 60
$ws.Range("G403").Value = 2232.6
$ws.Range("B407").Value = 48547.82
$ws.Range("F412").Value = 59
$ws.Range("G412").Value = 5796.75
$ws.Range("F416").Value = 67
$ws.Range("G416").Value = 2507.14
$ws.Range("F417").Value = 64
$ws.Range("G417").Value = 11932.16
$ws.Range("F418").Value = 91
$ws.Range("G418").Value = 3025.75
$ws.Range("B424").Value = 45548.43
$ws.Range("F458").Value = 1008
$ws.Range("G458").Value = 6632.64
$ws.Range("F463").Value = 301
$ws.Range("G463").Value = 4945.43
$ws.Range("B465").Value = 106907.38
$ws.Range("F470").Value = 50
$ws.Range("G470").Value = 1740.5
$ws.Range("B472").Value = 9393.1
$ws.Range("F508").Value = 651
$ws.Range("G508").Value = 4426.8
$ws.Range("F512").Value = 422
$ws.Range("G512").Value = 2827.4
$ws.Range("F514").Value = 139
$ws.Range("G514").Value = 3666.82
$ws.Range("B516").Value = 40861.76
$ws.Range("F525").Value = 82
$ws.Range("G525").Value = 5075.8
$ws.Range("B537").Value = 27452.7
$ws.Range("F574").Value = 6
$ws.Range("G574").Value = 1704.6
$ws.Range("F575").Value = 36
$ws.Range("G575").Value = 1005.48
$ws.Range("F576").Value = 0
$ws.Range("G576").Value = 0
$ws.Range("F577").Value = 121
$ws.Range("G577").Value = 3315.4
$ws.Range("F581").Value = 72
$ws.Range("G581").Value = 8843.76
$ws.Range("B583").Value = 41437.49
$ws.Range("F585").Value = 117
$ws.Range("G585").Value = 15274.35
$ws.Range("F589").Value = 86
$ws.Range("G589").Value = 2339.2
$ws.Range("F591").Value = 25
$ws.Range("G591").Value = 680
$ws.Range("B592").Value = 64292.09
$ws.Range("F614").Value = 177
$ws.Range("G614").Value = 7642.86
$ws.Range("F619").Value = 9
$ws.Range("G619").Value = 393.84
$ws.Range("B620").Value = 31384.09
$ws.Range("F669").Value = 0
$ws.Range("G669").Value = 0
$ws.Range("B676").Value = 19891.23
$ws.Range("F679").Value = 9
$ws.Range("G679").Value = 995.67
$ws.Range("F683").Value = 47
$ws.Range("G683").Value = 1948.15
$ws.Range("B688").Value = 43736.43
$ws.Range("F699").Value = 54
$ws.Range("G699").Value = 5557.14
$ws.Range("B706").Value = 23017.05
$ws.Range("F709").Value = 46
$ws.Range("G709").Value = 1720.4
$ws.Range("F711").Value = 98
$ws.Range("G711").Value = 3665.2
$ws.Range("B714").Value = 9300.700000000001
$ws.Range("F753").Value = 6
$ws.Range("G753").Value = 197.28
$ws.Range("F756").Value = 5
$ws.Range("G756").Value = 149.6
$ws.Range("B757").Value = 2391.59
$ws.Range("F759").Value = 1004
$ws.Range("G759").Value = 163762.44
$ws.Range("F761").Value = 212
$ws.Range("G761").Value = 16353.68
$ws.Range("F762").Value = 43
$ws.Range("G762").Value = 6357.98
$ws.Range("B764").Value = 193452.34
$ws.Range("B770").Value = 2373661.29
$ws.Range("B771").Value = 2373661.29
